# Insert a new data row at row 23 (pushing the existing rows 23-86 down to
# 24-87, dimension grows from A1:R86 to A1:R87), and populate the new row
# with a new "Acelga" price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23..86 down to 24..87, leaving a blank row at 23.
$ws.Rows.Item(23).Insert()

# Fill the newly inserted row 23 with the new record.
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44987
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100112009
$ws.Range("G23").Value = "Acelga"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 2500
$ws.Range("M23").Value = 2250
$ws.Range("N23").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 750
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = "Hortaliza"
